$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (e.g. "10.90"). Excel auto-coerces a
# numeric-looking string assigned via .Value into a real Number, so any cell
# whose new price looks like a plain number is switched to Text format first
# -- this keeps the stored cell value a string, matching the source data.

$ws.Range("D2").Value = "67.108.11"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "2.468.25"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.59"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.93"
$ws.Range("E6").Value = "  +2.26%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.94"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "2.922.82"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.35"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").Value = "67.025.19"

$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "2.425.31"
$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "7.45"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "10.87"
$ws.Range("E19").Value = "  -2.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.35"
$ws.Range("E20").Value = "  -1.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.02"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.33"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.19"
$ws.Range("E24").Value = "  -1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.17"
$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("D27").Value = "2.594.36"
$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").Value = "0.0₃0899"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "498.69"
$ws.Range("E30").Value = "  -4.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.70"
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("E33").Value = "  -1.49%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.119"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.87"
$ws.Range("E36").Value = "  +2.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.67"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.09"
$ws.Range("E38").Value = "  -1.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.32"
$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.68"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.326"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.81"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "142.09"
$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.508"
$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("D48").Value = "0.0₆0252"
$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.56"
$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.581"
$ws.Range("E51").Value = "  -0.21%  "
